$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B ---
$ws.Cells.Item(2, 2).Value = "NSE:KIMS"
$ws.Cells.Item(3, 2).Value = ""
$ws.Cells.Item(4, 2).Value = ""
$ws.Cells.Item(5, 2).Value = ""
$ws.Cells.Item(6, 2).Value = ""
$ws.Cells.Item(7, 2).Value = ""
$ws.Cells.Item(8, 2).Value = ""
$ws.Cells.Item(9, 2).Value = ""

# --- Update column E (row 2 only) ---
$ws.Cells.Item(2, 5).Value = "NSE:IDFC"

# --- Extend column A style down to row 69 (copy style from A35) before setting new values ---
$ws.Range("A35").Copy($ws.Range("A36:A69"))

# --- Update column A sequential numbers (0-based) rows 2-69 ---
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(51, 1).Value = 49
$ws.Cells.Item(52, 1).Value = 50
$ws.Cells.Item(53, 1).Value = 51
$ws.Cells.Item(54, 1).Value = 52
$ws.Cells.Item(55, 1).Value = 53
$ws.Cells.Item(56, 1).Value = 54
$ws.Cells.Item(57, 1).Value = 55
$ws.Cells.Item(58, 1).Value = 56
$ws.Cells.Item(59, 1).Value = 57
$ws.Cells.Item(60, 1).Value = 58
$ws.Cells.Item(61, 1).Value = 59
$ws.Cells.Item(62, 1).Value = 60
$ws.Cells.Item(63, 1).Value = 61
$ws.Cells.Item(64, 1).Value = 62
$ws.Cells.Item(65, 1).Value = 63
$ws.Cells.Item(66, 1).Value = 64
$ws.Cells.Item(67, 1).Value = 65
$ws.Cells.Item(68, 1).Value = 66
$ws.Cells.Item(69, 1).Value = 67

# --- Update column C for all rows 2-69 ---
$ws.Cells.Item(2, 3).Value = "NSE:ABBOTINDIA"
$ws.Cells.Item(3, 3).Value = "NSE:ADORWELD"
$ws.Cells.Item(4, 3).Value = "NSE:AGI"
$ws.Cells.Item(5, 3).Value = "NSE:AMRUTANJAN"
$ws.Cells.Item(6, 3).Value = "NSE:ASALCBR"
$ws.Cells.Item(7, 3).Value = "NSE:ATULAUTO"
$ws.Cells.Item(8, 3).Value = "NSE:AXSENSEX"
$ws.Cells.Item(9, 3).Value = "NSE:BANARISUG"
$ws.Cells.Item(10, 3).Value = "NSE:BASML"
$ws.Cells.Item(11, 3).Value = "NSE:BHARATGEAR"
$ws.Cells.Item(12, 3).Value = "NSE:BOSCHLTD"
$ws.Cells.Item(13, 3).Value = "NSE:DEEPENR"
$ws.Cells.Item(14, 3).Value = "NSE:DHAMPURSUG"
$ws.Cells.Item(15, 3).Value = "NSE:ELECON"
$ws.Cells.Item(16, 3).Value = "NSE:EMAMIREAL"
$ws.Cells.Item(17, 3).Value = "NSE:FOSECOIND"
$ws.Cells.Item(18, 3).Value = "NSE:FSL"
$ws.Cells.Item(19, 3).Value = "NSE:GENESYS"
$ws.Cells.Item(20, 3).Value = "NSE:GIPCL"
$ws.Cells.Item(21, 3).Value = "NSE:GMMPFAUDLR"
$ws.Cells.Item(22, 3).Value = "NSE:GOCLCORP"
$ws.Cells.Item(23, 3).Value = "NSE:GREAVESCOT"
$ws.Cells.Item(24, 3).Value = "NSE:HBSL"
$ws.Cells.Item(25, 3).Value = "NSE:HDFCSENSEX"
$ws.Cells.Item(26, 3).Value = "NSE:HESTERBIO"
$ws.Cells.Item(27, 3).Value = "NSE:HFCL"
$ws.Cells.Item(28, 3).Value = "NSE:HISARMETAL"
$ws.Cells.Item(29, 3).Value = "NSE:HONAUT"
$ws.Cells.Item(30, 3).Value = "NSE:IEL"
$ws.Cells.Item(31, 3).Value = "NSE:ISMTLTD"
$ws.Cells.Item(32, 3).Value = "NSE:JAGSNPHARM"
$ws.Cells.Item(33, 3).Value = "NSE:JAICORPLTD"
$ws.Cells.Item(34, 3).Value = "NSE:JINDWORLD"
$ws.Cells.Item(35, 3).Value = "NSE:KBCGLOBAL"
$ws.Cells.Item(36, 3).Value = "NSE:KECL"
$ws.Cells.Item(37, 3).Value = "NSE:KESORAMIND"
$ws.Cells.Item(38, 3).Value = "NSE:KHAICHEM"
$ws.Cells.Item(39, 3).Value = "NSE:KOTHARIPRO"
$ws.Cells.Item(40, 3).Value = "NSE:KREBSBIO"
$ws.Cells.Item(41, 3).Value = "NSE:LAXMICOT"
$ws.Cells.Item(42, 3).Value = "NSE:LLOYDSENGG"
$ws.Cells.Item(43, 3).Value = "NSE:MAHABANK"
$ws.Cells.Item(44, 3).Value = "NSE:MAXIND"
$ws.Cells.Item(45, 3).Value = "NSE:MGEL"
$ws.Cells.Item(46, 3).Value = "NSE:MRF"
$ws.Cells.Item(47, 3).Value = "NSE:MUKANDLTD"
$ws.Cells.Item(48, 3).Value = "NSE:NAGREEKEXP"
$ws.Cells.Item(49, 3).Value = "NSE:NAVINIFTY"
$ws.Cells.Item(50, 3).Value = "NSE:NETWORK18"
$ws.Cells.Item(51, 3).Value = "NSE:NFL"
$ws.Cells.Item(52, 3).Value = "NSE:NIACL"
$ws.Cells.Item(53, 3).Value = "NSE:NIPPOBATRY"
$ws.Cells.Item(54, 3).Value = "NSE:NRBBEARING"
$ws.Cells.Item(55, 3).Value = "NSE:NTPC"
$ws.Cells.Item(56, 3).Value = "NSE:ONWARDTEC"
$ws.Cells.Item(57, 3).Value = "NSE:PAISALO"
$ws.Cells.Item(58, 3).Value = "NSE:PGIL"
$ws.Cells.Item(59, 3).Value = "NSE:PLASTIBLEN"
$ws.Cells.Item(60, 3).Value = "NSE:POWERGRID"
$ws.Cells.Item(61, 3).Value = "NSE:PREMEXPLN"
$ws.Cells.Item(62, 3).Value = "NSE:PRIVISCL"
$ws.Cells.Item(63, 3).Value = "NSE:RAMRAT"
$ws.Cells.Item(64, 3).Value = "NSE:RBLBANK"
$ws.Cells.Item(65, 3).Value = "NSE:RENUKA"
$ws.Cells.Item(66, 3).Value = "NSE:RESPONIND"
$ws.Cells.Item(67, 3).Value = "NSE:RPGLIFE"
$ws.Cells.Item(68, 3).Value = "NSE:RSWM"
$ws.Cells.Item(69, 3).Value = "NSE:SAKHTISUG"
